# Reproduce the commit "Added All Java and Excel File":
#  - rename "Sheet2" -> "Sheet1"
#  - rename "Sheet3" -> "EmpDetails" and give it an "Address" header in B1
#  - append three brand-new sheets "EmpDetails1" / "EmpDetails2" / "EmpDetails3",
#    each carrying a "Name" header in A1
#  - leave "EmpDetails1" as the active/selected tab (B1 selected)
#  - move the selection on "WorkFlowData" from B5 to D7 (and drop its tab focus)

$wb = $excel.ActiveWorkbook

# --- rename the existing blank sheets -------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet2")
$sheet1.Name = "Sheet1"

$empDetails = $wb.Worksheets.Item("Sheet3")
$empDetails.Name = "EmpDetails"

# --- append the three new "EmpDetails" sheets after EmpDetails -----------
# Write the "Name" header cells first so the new shared-string entries are
# created in the same order as the target workbook (Name before Address).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$empDetails1 = $wb.Worksheets.Add($null, $lastSheet)
$empDetails1.Name = "EmpDetails1"
$empDetails1.Range("A1").Value = "Name"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$empDetails2 = $wb.Worksheets.Add($null, $lastSheet)
$empDetails2.Name = "EmpDetails2"
$empDetails2.Range("A1").Value = "Name"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$empDetails3 = $wb.Worksheets.Add($null, $lastSheet)
$empDetails3.Name = "EmpDetails3"
$empDetails3.Range("A1").Value = "Name"

# --- give EmpDetails its "Address" header ---------------------------------
$empDetails.Range("B1").Value = "Address"

# --- selections: WorkFlowData moves to D7, EmpDetails1 ends up active ----
$workFlowData = $wb.Worksheets.Item("WorkFlowData")
[void]$workFlowData.Range("D7").Select()

[void]$empDetails1.Activate()
[void]$empDetails1.Range("B1").Select()
